# This script inserts two new data rows (rows 828 and 829) into the single
# worksheet of the workbook, shifting all existing rows from 828 onward
# down by two positions (so old row 828 becomes row 830, etc.), and then
# populates the two newly inserted rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 828, shifting everything below down.
$ws.Rows.Item(828).Resize(2).Insert(-4121)

# --- New row 828 ---
$ws.Cells.Item(828, 1).Value  = 4
$ws.Cells.Item(828, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(828, 3).Value  = "Los Lagos"
$ws.Cells.Item(828, 4).Value  = 45154
$ws.Cells.Item(828, 5).Value  = 10
$ws.Cells.Item(828, 6).Value  = "Fruta"
$ws.Cells.Item(828, 7).Value  = 100108
$ws.Cells.Item(828, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(828, 9).Value  = 100108006
$ws.Cells.Item(828, 10).Value = "Plátano"
$ws.Cells.Item(828, 11).Value = "Sin especificar"
$ws.Cells.Item(828, 12).Value = "Pintón"
$ws.Cells.Item(828, 13).Value = 500
$ws.Cells.Item(828, 14).Value = 17000
$ws.Cells.Item(828, 15).Value = 17000
$ws.Cells.Item(828, 16).Value = 17000
$ws.Cells.Item(828, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(828, 18).Value = "Ecuador"
$ws.Cells.Item(828, 19).Value = 850
$ws.Cells.Item(828, 20).Value = 20

# --- New row 829 ---
$ws.Cells.Item(829, 1).Value  = 4
$ws.Cells.Item(829, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(829, 3).Value  = "Los Lagos"
$ws.Cells.Item(829, 4).Value  = 45154
$ws.Cells.Item(829, 5).Value  = 10
$ws.Cells.Item(829, 6).Value  = "Fruta"
$ws.Cells.Item(829, 7).Value  = 100108
$ws.Cells.Item(829, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(829, 9).Value  = 100108006
$ws.Cells.Item(829, 10).Value = "Plátano"
$ws.Cells.Item(829, 11).Value = "Sin especificar"
$ws.Cells.Item(829, 12).Value = "Primera Pintón"
$ws.Cells.Item(829, 13).Value = 800
$ws.Cells.Item(829, 14).Value = 19000
$ws.Cells.Item(829, 15).Value = 20000
$ws.Cells.Item(829, 16).Value = 19500
$ws.Cells.Item(829, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(829, 18).Value = "Ecuador"
$ws.Cells.Item(829, 19).Value = 975
$ws.Cells.Item(829, 20).Value = 20

# Make sure column D keeps the date number format used by the rest of the
# column (it is applied automatically by Insert, but set it explicitly too
# to be safe).
$ws.Range("D828:D829").NumberFormat = $ws.Range("D827").NumberFormat
